$wb = $excel.ActiveWorkbook

# --- Sheet references (by fixed index; tab order is unchanged by this edit) ---
$wsData  = $wb.Worksheets.Item(1)   # dss_detailed_timing_data
$wsSubs  = $wb.Worksheets.Item(2)   # dss_subroutines
$wsSrc   = $wb.Worksheets.Item(3)   # dss_source_code

# --- 1. dss_source_code: add the source-file caption (must be written before
#        the dss_detailed_timing_data caption so the shared-string table gets
#        the same ordering as the target file: idx 9 then idx 10) ---
$wsSrc.Range("A1").Value = "Source code filepath: LinDistFlow/powerflowpy/powerflowpy/dss_solve_detailed_timing.py"

# move the picture down two rows to make room for the new text above it
$pic = $wsSrc.Shapes.Item(1)
$pic.Width = 987.954251968504
$pic.Top = 32

# --- 2. dss_detailed_timing_data: insert a new title row above the old one ---
$wsData.Rows.Item(1).Insert()
$wsData.Range("A1").Value = "Source code: LinDist3Flow/powerflowpy/dss_timer.py"

# reset its selection back to the default (A1)
$wsData.Range("A1").Select()

# --- 3. dss_subroutines: sort the summary table descending by run time ---
$wsSubs.Sort.SortFields.Clear()
$wsSubs.Sort.SortFields.Add($wsSubs.Range("D48:D52"), 0, 2) | Out-Null
$wsSubs.Sort.SetRange($wsSubs.Range("C48:D52"))
$wsSubs.Sort.Header = 0
$wsSubs.Sort.Apply()

# widen the label column so the sorted labels are readable
$wsSubs.Columns.Item(3).ColumnWidth = 19.83

# --- 4. View state: dss_source_code keeps a non-default selection, then
#        dss_subroutines becomes the active / visible tab ---
$wsSrc.Range("V24").Select()
$wsSubs.Range("G51").Select()
$wsSubs.Activate()
